# Applies updated transition-probability values to Sheet1 of the workbook.
# These values come from re-running the simulation with more games / updated
# simulate-game + optimization logic, which shifted several matrix entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.2266666666666667
$ws.Range("C2").Value2 = 0.5022222222222222
$ws.Range("P2").Value2 = 0.1644444444444444
$ws.Range("S2").Value2 = 0.1066666666666667
$ws.Range("B3").Value2 = 0.01785714285714286
$ws.Range("C3").Value2 = 0.008928571428571428
$ws.Range("J3").Value2 = 0.008928571428571428
$ws.Range("P3").Value2 = 0.7767857142857143
$ws.Range("S3").Value2 = 0.1875
$ws.Range("J4").Value2 = 0.1153846153846154
$ws.Range("P4").Value2 = 0.5
$ws.Range("S4").Value2 = 0.3846153846153846
$ws.Range("B6").Value2 = 0.0546218487394958
$ws.Range("F6").Value2 = 0.07983193277310924
$ws.Range("J6").Value2 = 0.2478991596638656
$ws.Range("O6").Value2 = 0.01260504201680672
$ws.Range("Q6").Value2 = 0.1176470588235294
$ws.Range("R6").Value2 = 0.07983193277310924
$ws.Range("S6").Value2 = 0.407563025210084
$ws.Range("B7").Value2 = 0.06862745098039216
$ws.Range("D7").Value2 = 0.0196078431372549
$ws.Range("E7").Value2 = 0.004901960784313725
$ws.Range("F7").Value2 = 0.06372549019607843
$ws.Range("J7").Value2 = 0.1127450980392157
$ws.Range("O7").Value2 = 0.02450980392156863
$ws.Range("Q7").Value2 = 0.142156862745098
$ws.Range("S7").Value2 = 0.4803921568627451
$ws.Range("B8").Value2 = 0.08148148148148149
$ws.Range("D8").Value2 = 0.01234567901234568
$ws.Range("F8").Value2 = 0.05925925925925926
$ws.Range("J8").Value2 = 0.1037037037037037
$ws.Range("O8").Value2 = 0.01975308641975309
$ws.Range("Q8").Value2 = 0.1506172839506173
$ws.Range("R8").Value2 = 0.09135802469135802
$ws.Range("S8").Value2 = 0.4814814814814815
$ws.Range("B9").Value2 = 0.07407407407407407
$ws.Range("D9").Value2 = 0.01481481481481482
$ws.Range("F9").Value2 = 0.05925925925925926
$ws.Range("J9").Value2 = 0.08888888888888889
$ws.Range("O9").Value2 = 0.01851851851851852
$ws.Range("Q9").Value2 = 0.1333333333333333
$ws.Range("R9").Value2 = 0.1185185185185185
$ws.Range("S9").Value2 = 0.4925925925925926
$ws.Range("B10").Value2 = 0.07247557003257329
$ws.Range("D10").Value2 = 0.01465798045602606
$ws.Range("F10").Value2 = 0.07899022801302931
$ws.Range("J10").Value2 = 0.09039087947882736
$ws.Range("O10").Value2 = 0.01791530944625407
$ws.Range("Q10").Value2 = 0.1767100977198697
$ws.Range("R10").Value2 = 0.0993485342019544
$ws.Range("S10").Value2 = 0.4495114006514658
$ws.Range("G11").Value2 = 0.1408450704225352
$ws.Range("J11").Value2 = 0.09154929577464789
$ws.Range("K11").Value2 = 0.2147887323943662
$ws.Range("L11").Value2 = 0.5316901408450704
$ws.Range("S11").Value2 = 0.02112676056338028
$ws.Range("G12").Value2 = 0.808641975308642
$ws.Range("J12").Value2 = 0.1172839506172839
$ws.Range("K12").Value2 = 0.006172839506172839
$ws.Range("L12").Value2 = 0.04320987654320987
$ws.Range("S12").Value2 = 0.02469135802469136
$ws.Range("G13").Value2 = 0.7777777777777778
$ws.Range("J13").Value2 = 0.2222222222222222
$ws.Range("G14").Value2 = 0.4
$ws.Range("J14").Value2 = 0.4
$ws.Range("S14").Value2 = 0.2
$ws.Range("F15").Value2 = 0.02512562814070352
$ws.Range("H15").Value2 = 0.1005025125628141
$ws.Range("I15").Value2 = 0.1256281407035176
$ws.Range("J15").Value2 = 0.3366834170854272
$ws.Range("K15").Value2 = 0.07035175879396985
$ws.Range("M15").Value2 = 0.01005025125628141
$ws.Range("N15").Value2 = 0.005025125628140704
$ws.Range("O15").Value2 = 0.04020100502512563
$ws.Range("S15").Value2 = 0.2864321608040201
$ws.Range("F16").Value2 = 0.01492537313432836
$ws.Range("H16").Value2 = 0.1492537313432836
$ws.Range("I16").Value2 = 0.1119402985074627
$ws.Range("J16").Value2 = 0.4626865671641791
$ws.Range("K16").Value2 = 0.08208955223880597
$ws.Range("M16").Value2 = 0.007462686567164179
$ws.Range("O16").Value2 = 0.03731343283582089
$ws.Range("S16").Value2 = 0.1343283582089552
$ws.Range("F17").Value2 = 0.01907356948228883
$ws.Range("H17").Value2 = 0.1689373297002725
$ws.Range("I17").Value2 = 0.1198910081743869
$ws.Range("J17").Value2 = 0.4414168937329701
$ws.Range("K17").Value2 = 0.05177111716621254
$ws.Range("M17").Value2 = 0.01634877384196185
$ws.Range("O17").Value2 = 0.04087193460490463
$ws.Range("S17").Value2 = 0.1416893732970027
$ws.Range("F18").Value2 = 0.02183406113537118
$ws.Range("H18").Value2 = 0.1222707423580786
$ws.Range("I18").Value2 = 0.1441048034934498
$ws.Range("J18").Value2 = 0.462882096069869
$ws.Range("K18").Value2 = 0.05240174672489083
$ws.Range("M18").Value2 = 0.01310043668122271
$ws.Range("O18").Value2 = 0.03493449781659388
$ws.Range("S18").Value2 = 0.148471615720524
$ws.Range("F19").Value2 = 0.01711156741957563
$ws.Range("H19").Value2 = 0.1909650924024641
$ws.Range("I19").Value2 = 0.1054072553045859
$ws.Range("J19").Value2 = 0.3620807665982204
$ws.Range("K19").Value2 = 0.1115674195756331
$ws.Range("M19").Value2 = 0.02190280629705681
$ws.Range("N19").Value2 = 0.002737850787132101
$ws.Range("O19").Value2 = 0.06433949349760439
$ws.Range("S19").Value2 = 0.1238877481177276

Write-Host "Updated transition matrix values on sheet $($ws.Name)"
